$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "Senior Manager, Data Science" -> "Senior Data Scientist"
# (title text in the Professional Work Experience table; the whole
#  phrase keeps its original bold / sz17 / szCs17 run formatting)
# ---------------------------------------------------------------------
$find1 = $d.Content
$find1.Find.ClearFormatting()
$find1.Find.Execute("Senior Manager, Data Science", $false, $false, $false, $false, $false, `
                     $true, 1, $false, "Senior Data Scientist", 2) | Out-Null
Write-Output "Updated title to 'Senior Data Scientist'"

# ---------------------------------------------------------------------
# Change 2: append ", and Tableau" right after "using Python and SQL"
# (before the trailing period). "Tableau" is bold; the connector text
# ", and " is not.
# ---------------------------------------------------------------------
$find2 = $d.Content
$find2.Find.ClearFormatting()
$find2.Find.Execute("using Python and SQL", $false, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0) | Out-Null

# Collapse to the end of the found range, right before the trailing "."
$insertPoint = $find2.Duplicate
$insertPoint.Collapse(0)   # wdCollapseEnd

$insertPoint.InsertAfter(", and ")
$insertPoint.Font.Bold = $false
$insertPoint.LanguageID = 1066   # wdVietnamese
Write-Output "Inserted connector text ', and '"

$tableauPoint = $insertPoint.Duplicate
$tableauPoint.Collapse(0)  # wdCollapseEnd

$tableauPoint.InsertAfter("Tableau")
$tableauPoint.Font.Bold = $true
$tableauPoint.LanguageID = 1066   # wdVietnamese
Write-Output "Inserted 'Tableau'"
